# Modulos_Porta-Internol-v2.xlsx
# "Se actualiza listado de módulos del portal interno"
#
# Renames the two module-name labels that appear (each) in three places
# on the sheet, and brings the indentation/style of a couple of cells in
# line with their sibling rows. Also nudges the view back to a plain,
# zoomed-out state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the module labels.
#    "ltl_disponibles(dentro del módulo)" -> "Ltl"
#    "cd_disponibles(dentro del módulo)"  -> "CrossDock"
#    Each string shows up in three different cells; updating the cell
#    text in place rewrites the shared-string entry for every occurrence.
# ---------------------------------------------------------------------
$ws.Range("H7").Value  = "Ltl"
$ws.Range("H8").Value  = "CrossDock"
$ws.Range("H10").Value = "Ltl"
$ws.Range("H11").Value = "CrossDock"
$ws.Range("F23").Value = "Ltl"
$ws.Range("F24").Value = "CrossDock"

# ---------------------------------------------------------------------
# 2. Style touch-ups.
#    A4/A5 pick up the same left-indent style already used by the rest
#    of their row (B4/B5). F24 picks up the deeper-indent style already
#    used by the "CrossDock" entries in column H (e.g. H8).
#    Copy/PasteSpecial(formats) carries over the style without touching
#    the destination cell's value.
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Range("B4").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)

$ws.Range("B5").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)

$ws.Range("H8").Copy()
$ws.Range("F24").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. View: drop the saved scroll position / selection and zoom out to
#    85% as the reopened sheet did.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
$excel.ActiveWindow.Zoom = 85
